$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-12-24 20:49:47"

# Update the timestamp in column O for every data row (rows 2 through 410)
for ($row = 2; $row -le 410; $row++) {
    $ws.Range("O$row").Value = $newTimestamp
}

# Row 233 also had its productAriaLabel text (column M) updated
$ws.Range("M233").Value = "Betty Bossi Kuchenteig -35% Fett - Online kein Bestand 2.10 Schweizer Franken"
